$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain text so numeric-looking values
# ("218.01", "1.00", etc.) are not auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '28.338.95'
$ws.Range("E2").Value = '  +4.21%  '

$ws.Range("D3").Value = '1.710.26'
$ws.Range("E3").Value = '  +1.42%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = '218.01'
$ws.Range("E5").Value = '  +0.87%  '

$ws.Range("D6").Value = '0.522'
$ws.Range("E6").Value = '  +0.33%  '

$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("D8").Value = '24.25'
$ws.Range("E8").Value = '  +4.62%  '

$ws.Range("E9").Value = '  +2.28%  '

$ws.Range("D10").Value = '0.0631'
$ws.Range("E10").Value = '  +0.76%  '

$ws.Range("E11").Value = '  -0.16%  '

$ws.Range("D12").Value = '1.954.50'
$ws.Range("E12").Value = '  +1.54%  '

$ws.Range("D13").Value = '1.711.92'
$ws.Range("E13").Value = '  +1.40%  '

$ws.Range("E14").Value = '  +0.00%  '

$ws.Range("E15").Value = '  +0.80%  '

$ws.Range("D16").Value = '67.31'
$ws.Range("E16").Value = '  +0.55%  '

$ws.Range("D17").Value = '253.41'
$ws.Range("E17").Value = '  +7.27%  '

$ws.Range("D18").Value = '28.302.52'
$ws.Range("E18").Value = '  +4.08%  '

$ws.Range("E19").Value = '  +0.35%  '

$ws.Range("D20").Value = '7.73'
$ws.Range("E20").Value = '  -3.75%  '

$ws.Range("E21").Value = '  -0.25%  '

$ws.Range("D22").Value = '4.56'
$ws.Range("E22").Value = '  -0.17%  '

$ws.Range("E23").Value = '  -0.32%  '

$ws.Range("E24").Value = '  -1.82%  '

$ws.Range("D25").Value = '148.06'
$ws.Range("E25").Value = '  +0.35%  '

$ws.Range("D26").Value = '7.36'
$ws.Range("E26").Value = '  +0.39%  '

$ws.Range("D27").Value = '16.60'
$ws.Range("E27").Value = '  +0.90%  '

$ws.Range("E28").Value = '  +0.33%  '

$ws.Range("E29").Value = '  +0.07%  '

$ws.Range("E30").Value = '  +0.93%  '

$ws.Range("E31").Value = '  +2.84%  '

$ws.Range("D32").Value = '3.39'
$ws.Range("E32").Value = '  +0.49%  '

$ws.Range("D33").Value = '1.481.70'
$ws.Range("E33").Value = '  -3.98%  '

$ws.Range("E34").Value = '  -1.54%  '

$ws.Range("E35").Value = '  -2.29%  '

$ws.Range("D36").Value = '0.963'
$ws.Range("E36").Value = '  +1.74%  '

$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("E38").Value = '  -1.57%  '

$ws.Range("E39").Value = '  +0.35%  '

$ws.Range("D40").Value = '1.04'
$ws.Range("E40").Value = '  -1.53%  '

$ws.Range("D41").Value = '69.49'
$ws.Range("E41").Value = '  +0.52%  '

$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("D43").Value = '5.64'
$ws.Range("E43").Value = '  -2.05%  '

$ws.Range("D44").Value = '1.859.38'
$ws.Range("E44").Value = '  +1.40%  '

$ws.Range("E45").Value = '  +0.09%  '

$ws.Range("D46").Value = '0.800'
$ws.Range("E46").Value = '  +1.31%  '

$ws.Range("E47").Value = '  +7.00%  '

$ws.Range("D48").Value = '89.93'
$ws.Range("E48").Value = '  -0.32%  '

$ws.Range("D49").Value = '0.0₆0108'
$ws.Range("E49").Value = '  -4.08%  '

$ws.Range("E50").Value = '  -0.68%  '

$ws.Range("D51").Value = '8.02'
$ws.Range("E51").Value = '  -3.21%  '

# Restore the default (unstyled) cell style now that the text values are set,
# so no extra explicit formatting is left behind on the Price column.
$priceRange.Style = "Normal"
